$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New device_master_h rows (added two new Mac-Address batches: "30" and "31")
$rows = @(
    @{ Row = 147; Id = 3000166; Name = "Finger Print Scanner 30"; Mac = "D6-15-AC-80-6B-86"; Serial = "BS563Q2230814"; DspecId = 165 },
    @{ Row = 148; Id = 3000167; Name = "IRIS Scanner 30";         Mac = "6D-58-E2-DF-74-34"; Serial = "BS563Q2230815"; DspecId = 327 },
    @{ Row = 149; Id = 3000168; Name = "Web Camera 30";           Mac = "E2-A8-56-86-15-30"; Serial = "BS563Q2230816"; DspecId = 736 },
    @{ Row = 150; Id = 3000169; Name = "Document Scanner 30";     Mac = "72-E8-B9-FD-63-65"; Serial = "BS563Q2230817"; DspecId = 801 },
    @{ Row = 151; Id = 3000170; Name = "Printer 30";               Mac = "D3-F3-A4-50-AD-12"; Serial = "BS563Q2230818"; DspecId = 920 },
    @{ Row = 152; Id = 3000171; Name = "Finger Print Scanner 31"; Mac = "06-16-D0-0B-A6-E4"; Serial = "BS563Q2230819"; DspecId = 165 },
    @{ Row = 153; Id = 3000172; Name = "IRIS Scanner 31";         Mac = "21-78-45-AC-E9-20"; Serial = "BS563Q2230820"; DspecId = 327 },
    @{ Row = 154; Id = 3000173; Name = "Web Camera 31";           Mac = "3C-E8-87-99-DB-FA"; Serial = "BS563Q2230821"; DspecId = 736 },
    @{ Row = 155; Id = 3000174; Name = "Document Scanner 31";     Mac = "BF-55-53-98-40-08"; Serial = "BS563Q2230822"; DspecId = 801 },
    @{ Row = 156; Id = 3000175; Name = "Printer 31";               Mac = "5A-43-36-46-22-EB"; Serial = "BS563Q2230823"; DspecId = 920 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $r.Id
    $ws.Range("B$rowNum").Value = $r.Name
    $ws.Range("C$rowNum").Value = $r.Mac
    $ws.Range("D$rowNum").Value = $r.Serial
    $ws.Range("F$rowNum").Value = $r.DspecId
    $ws.Range("G$rowNum").Value = "eng"
    $ws.Range("H$rowNum").Value = $true
    $ws.Range("H$rowNum").HorizontalAlignment = -4131
    $ws.Range("I$rowNum").Value = "superadmin"
    $ws.Range("J$rowNum").Value = "now()"
    $ws.Range("K$rowNum").Value = "now()"
}

# Match the saved view state recorded in the diff: scrolled down, selection on D145
$excel.ActiveWindow.ScrollRow = 139
$ws.Range("D145").Select() | Out-Null
